$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(152).Insert()

$ws.Cells.Item(152, 1).Value = 5
$ws.Cells.Item(152, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(152, 3).Value = "Maule"
$ws.Cells.Item(152, 4).Value = 44741
$ws.Cells.Item(152, 5).Value = 7
$ws.Cells.Item(152, 6).Value = "Fruta"
$ws.Cells.Item(152, 7).Value = 100108
$ws.Cells.Item(152, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(152, 9).Value = 100108005
$ws.Cells.Item(152, 10).Value = "Piña"
$ws.Cells.Item(152, 11).Value = "Caramelo"
$ws.Cells.Item(152, 12).Value = "Segunda"
$ws.Cells.Item(152, 13).Value = 250
$ws.Cells.Item(152, 14).Value = 17000
$ws.Cells.Item(152, 15).Value = 17000
$ws.Cells.Item(152, 16).Value = 17000
$ws.Cells.Item(152, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(152, 18).Value = "Ecuador"
$ws.Cells.Item(152, 19).Value = 1214
$ws.Cells.Item(152, 20).Value = 14
